$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) stores numeric-looking values (e.g. "62.279.07",
# "1.00") as plain text in the original workbook. Assigning such strings
# straight to .Value lets Excel coerce the unambiguous ones (single decimal
# point) into real numbers, so for each of those cells we briefly force a
# text number format while writing the value, then restore the default
# "Normal" cell style so the cell ends up unstyled again, matching the
# original file (NumberFormat must be set per cell - applying it to a
# multi-area union range only affects the first area).
$dRefs = @('D2', 'D3', 'D5', 'D6', 'D8', 'D9', 'D12', 'D13', 'D14', 'D17', 'D18', 'D19', 'D20', 'D22', 'D24', 'D26', 'D27', 'D30', 'D31', 'D32', 'D34', 'D37', 'D39', 'D40', 'D41', 'D43', 'D44', 'D45', 'D46', 'D48', 'D49', 'D50')
foreach ($ref in $dRefs) { $ws.Range($ref).NumberFormat = "@" }

$ws.Range('D2').Value = '62.279.07'
$ws.Range('E2').Value = '  +0.19%  '
$ws.Range('D3').Value = '2.430.96'
$ws.Range('E3').Value = '  +0.75%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').Value = '563.32'
$ws.Range('E5').Value = '  +0.41%  '
$ws.Range('D6').Value = '144.61'
$ws.Range('E6').Value = '  +1.19%  '
$ws.Range('E7').Value = '  -0.11%  '
$ws.Range('D8').Value = '0.532'
$ws.Range('E8').Value = '  +0.39%  '
$ws.Range('D9').Value = '2.430.14'
$ws.Range('E10').Value = '  +1.01%  '
$ws.Range('E11').Value = '  +0.29%  '
$ws.Range('D12').Value = '5.24'
$ws.Range('E12').Value = '  -1.40%  '
$ws.Range('D13').Value = '0.350'
$ws.Range('E13').Value = '  -0.15%  '
$ws.Range('D14').Value = '26.57'
$ws.Range('E14').Value = '  +3.64%  '
$ws.Range('E15').Value = '  +0.71%  '
$ws.Range('D17').Value = '62.075.19'
$ws.Range('E17').Value = '  -0.11%  '
$ws.Range('D18').Value = '2.424.21'
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('D19').Value = '11.22'
$ws.Range('E19').Value = '  -0.37%  '
$ws.Range('D20').Value = '323.95'
$ws.Range('E20').Value = '  +0.29%  '
$ws.Range('E21').Value = '  +1.11%  '
$ws.Range('D22').Value = '4.15'
$ws.Range('E22').Value = '  -0.45%  '
$ws.Range('E23').Value = '  -0.05%  '
$ws.Range('D24').Value = '67.40'
$ws.Range('E24').Value = '  +2.53%  '
$ws.Range('E25').Value = '  +1.19%  '
$ws.Range('D26').Value = '8.81'
$ws.Range('E26').Value = '  -2.26%  '
$ws.Range('D27').Value = '553.19'
$ws.Range('E27').Value = '  -3.82%  '
$ws.Range('E28').Value = '  +0.74%  '
$ws.Range('E29').Value = '  +0.85%  '
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  -0.22%  '
$ws.Range('D31').Value = '8.28'
$ws.Range('E31').Value = '  +1.01%  '
$ws.Range('D32').Value = '1.43'
$ws.Range('E32').Value = '  -0.81%  '
$ws.Range('E33').Value = '  -0.82%  '
$ws.Range('D34').Value = '1.88'
$ws.Range('E34').Value = '  +0.51%  '
$ws.Range('E35').Value = '  -0.77%  '
$ws.Range('E36').Value = '  -0.14%  '
$ws.Range('D37').Value = '4.83'
$ws.Range('E37').Value = '  +2.36%  '
$ws.Range('E38').Value = '  +0.06%  '
$ws.Range('D39').Value = '5.57'
$ws.Range('E39').Value = '  -0.24%  '
$ws.Range('D40').Value = '18.69'
$ws.Range('E40').Value = '  +0.17%  '
$ws.Range('D41').Value = '149.97'
$ws.Range('E41').Value = '  -2.06%  '
$ws.Range('E42').Value = '  -0.62%  '
$ws.Range('D43').Value = '1.00'
$ws.Range('E43').Value = '  +0.04%  '
$ws.Range('D44').Value = '2.33'
$ws.Range('E44').Value = '  +1.70%  '
$ws.Range('D45').Value = '147.54'
$ws.Range('E45').Value = '  -0.76%  '
$ws.Range('D46').Value = '3.68'
$ws.Range('E46').Value = '  +1.01%  '
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('D48').Value = '20.24'
$ws.Range('E48').Value = '  +0.96%  '
$ws.Range('D49').Value = '0.597'
$ws.Range('E49').Value = '  +0.82%  '
$ws.Range('D50').Value = '0.0924'
$ws.Range('E50').Value = '  +0.87%  '
$ws.Range('E51').Value = '  +0.95%  '

foreach ($ref in $dRefs) { $ws.Range($ref).Style = "Normal" }
